$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-13 01:19:48"

for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
